# Apply edits to match the target diff: re-sort several rows (identity/score
# swaps between tied positions) and recompute the "matrices" ranking scores.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 13.4562512144795
$ws.Range("H3").Value = 13.03225951465968
$ws.Range("D4").Value = 21
$ws.Range("E4").Value = "5c0e89c6c323400001e6c4a5"
$ws.Range("F4").Value = "Bri"
$ws.Range("H4").Value = 8.49956431831203
$ws.Range("I4").Value = "Black or African American"
$ws.Range("H5").Value = 8.20087022370102
$ws.Range("D6").Value = 22
$ws.Range("E6").Value = "608b14a312c099ac00b721b6"
$ws.Range("F6").Value = "Khushi"
$ws.Range("H6").Value = 8.064168822103696
$ws.Range("I6").Value = "Asian"
$ws.Range("H7").Value = 5.333931338090698
$ws.Range("H8").Value = 5.316202313826643
$ws.Range("H9").Value = 5.294121455295787
$ws.Range("D10").Value = 34
$ws.Range("E10").Value = "5e96194b0a9fe909389e9f7b"
$ws.Range("F10").Value = "Tina"
$ws.Range("H10").Value = 4.302967855272668
$ws.Range("D11").Value = 35
$ws.Range("E11").Value = "6077db0613ce87b4a62a78f9"
$ws.Range("F11").Value = "Lori"
$ws.Range("H11").Value = 4.23075704731449
$ws.Range("H12").Value = 2.016984074606204
$ws.Range("H13").Value = 1.054877676087834
$ws.Range("H14").Value = 14.32124806351207
$ws.Range("H15").Value = 13.25929927578149
$ws.Range("H16").Value = 8.119568600285705
$ws.Range("H17").Value = 7.411316600731239
$ws.Range("H18").Value = 6.387209556654361
$ws.Range("H19").Value = 6.159915716578424
$ws.Range("D20").Value = 32
$ws.Range("E20").Value = "60bf9943e4e04642d4634ecc"
$ws.Range("F20").Value = "Jamarii"
$ws.Range("H20").Value = 5.496872041548905
$ws.Range("I20").Value = "Black or African American"
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = "60c2341fe95d71ee52c043f0"
$ws.Range("F21").Value = "Matthew"
$ws.Range("H21").Value = 5.269135301367183
$ws.Range("D22").Value = 33
$ws.Range("E22").Value = "60b322994d0b901954690036"
$ws.Range("F22").Value = "Brennan"
$ws.Range("H22").Value = 5.262741384947466
$ws.Range("I22").Value = "White"
$ws.Range("H23").Value = 3.214622224054206
$ws.Range("H24").Value = 1.200783564827204
$ws.Range("H25").Value = 0.09961791273931464
